# Aula 40 - Finalizando a acao de exclusao - testado e funcionando
#
# Adds 4 new rows (36-39) to the bottom of the "anotacoes" worksheet
# documenting aula 40 ("Finalizando a ação de exclusão").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 already carries the exact cell style pattern the new rows need
# (B/C -> style 5, D -> style 6, E -> style 1). Copy that 4-cell slice's
# formatting onto the new rows first so the new cells reuse the existing
# cellXfs entries instead of the engine minting new ones, and so B/C/D/E
# all stay inside a tight "2:5" row span just like every other data row.
$ws.Range("B31:E31").Copy() | Out-Null
$ws.Range("B36:E39").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Populate the brand-new text in the same order the authoring tool used,
# so new entries land in xl/sharedStrings.xml with the expected indices
# (65 already exists; new ones become 66..70 in this exact sequence):
#   66 - the "2:14 ... transformando o modal ..." note (row 37, col E)
#   67 - "40. Finalizando a ação de exclusão"          (col D, every row)
#   68 - the "0:48 ... primeiro uso de MODAL ..." note (row 36, col E)
#   69 - the "5:05 ... javascrpit ..." note             (row 38, col E)
#   70 - the "10:53 ... em resumo ..." note              (row 39, col E)
$ws.Cells.Item(37, 5).Value2 = "2:14`n8. Departamento: Controller & View`n40. Finalizando a ação de exclusão`ntransformando o modal em um fragment com th:fragment"
$ws.Cells.Item(36, 4).Value2 = "40. Finalizando a ação de exclusão"
$ws.Cells.Item(36, 5).Value2 = "0:48`nprimeiro uso de MODAL (um componente do bootstrap) no projeto. Uma espécie de janela de confirmação antes de excluir um departamento."
$ws.Cells.Item(38, 5).Value2 = "5:05`nprimeira abordagem de javascrpit no projeto, para trabalhar com o botao de excluir na lista de departamentos."
$ws.Cells.Item(39, 5).Value2 = "10:53`nem resumo, na aula 40 implementamos o botão de excluir no frontend, construimos um fragmento de pagina MODAL modelo de confirmação de exclusão antes de excluir o objeto da lista, onde a exclusão de fato do objeto acontece ao clicar em OK na mensagem, ou seja, a requisição URL de deletar acontecer ao CONFIRMAR no modal."

# Remaining cells that reuse already-known shared strings (no new indices).
$ws.Cells.Item(36, 2).Value2 = 40
$ws.Cells.Item(36, 3).Value2 = "8. Departamento: Controller & View"

$ws.Cells.Item(37, 2).Value2 = 40
$ws.Cells.Item(37, 3).Value2 = "8. Departamento: Controller & View"
$ws.Cells.Item(37, 4).Value2 = "40. Finalizando a ação de exclusão"

$ws.Cells.Item(38, 2).Value2 = 40
$ws.Cells.Item(38, 3).Value2 = "8. Departamento: Controller & View"
$ws.Cells.Item(38, 4).Value2 = "40. Finalizando a ação de exclusão"

$ws.Cells.Item(39, 2).Value2 = 40
$ws.Cells.Item(39, 3).Value2 = "8. Departamento: Controller & View"
$ws.Cells.Item(39, 4).Value2 = "40. Finalizando a ação de exclusão"

# Row heights as authored (wrapped text heights for each note's line count).
$ws.Rows.Item(36).RowHeight = 45
$ws.Rows.Item(37).RowHeight = 60
$ws.Rows.Item(38).RowHeight = 45
$ws.Rows.Item(39).RowHeight = 90

# Keep the view scrolled/selected the same way the authored workbook ended
# up after appending the new rows (active cell E39).
$ws.Application.Goto($ws.Range("A35"))
$ws.Range("E39").Select() | Out-Null
